$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# B6 / B7 change from numeric plant codes to text plant codes
$ws.Range("B6").Value = "A001"
$ws.Range("B7").Value = "Z001"

# A10:A15 become simple sequential line numbers, B10:B15 become a flat 555 each
$ws.Range("A10").Value = 1
$ws.Range("A11").Value = 2
$ws.Range("A12").Value = 3
$ws.Range("A13").Value = 4
$ws.Range("A14").Value = 5
$ws.Range("A15").Value = 6

$ws.Range("B10").Value = 555
$ws.Range("B11").Value = 555
$ws.Range("B12").Value = 555
$ws.Range("B13").Value = 555
$ws.Range("B14").Value = 555
$ws.Range("B15").Value = 555

# Page setup: A4 portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the visible selection to B10 (first row of the logistics table)
$ws.Range("B10").Select()

$wb.Save()
